$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d2 = $ws.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "42.750.66"
$d2.ClearFormats()
$ws.Range("E2").Value = "  -1.09%  "

$d3 = $ws.Range("D3")
$d3.NumberFormat = "@"
$d3.Value = "2.266.23"
$d3.ClearFormats()
$ws.Range("E3").Value = "  -1.63%  "

$ws.Range("E4").Value = "  -0.12%  "

$d5 = $ws.Range("D5")
$d5.NumberFormat = "@"
$d5.Value = "249.54"
$d5.ClearFormats()
$ws.Range("E5").Value = "  -0.85%  "

$d6 = $ws.Range("D6")
$d6.NumberFormat = "@"
$d6.Value = "0.632"
$d6.ClearFormats()
$ws.Range("E6").Value = "  -1.73%  "

$d7 = $ws.Range("D7")
$d7.NumberFormat = "@"
$d7.Value = "78.07"
$d7.ClearFormats()
$ws.Range("E7").Value = "  +6.01%  "

$ws.Range("E8").Value = "  -0.02%  "

$d9 = $ws.Range("D9")
$d9.NumberFormat = "@"
$d9.Value = "0.646"
$d9.ClearFormats()
$ws.Range("E9").Value = "  -2.90%  "

$d10 = $ws.Range("D10")
$d10.NumberFormat = "@"
$d10.Value = "40.43"
$d10.ClearFormats()
$ws.Range("E10").Value = "  +2.10%  "

$d11 = $ws.Range("D11")
$d11.NumberFormat = "@"
$d11.Value = "0.0962"
$d11.ClearFormats()
$ws.Range("E11").Value = "  -2.19%  "

$d12 = $ws.Range("D12")
$d12.NumberFormat = "@"
$d12.Value = "7.32"
$d12.ClearFormats()
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("E13").Value = "  -0.47%  "

$d14 = $ws.Range("D14")
$d14.NumberFormat = "@"
$d14.Value = "2.604.72"
$d14.ClearFormats()
$ws.Range("E14").Value = "  -1.59%  "

$d15 = $ws.Range("D15")
$d15.NumberFormat = "@"
$d15.Value = "15.03"
$d15.ClearFormats()
$ws.Range("E15").Value = "  -1.71%  "

$d16 = $ws.Range("D16")
$d16.NumberFormat = "@"
$d16.Value = "0.863"
$d16.ClearFormats()
$ws.Range("E16").Value = "  -4.28%  "

$d17 = $ws.Range("D17")
$d17.NumberFormat = "@"
$d17.Value = "2.273.17"
$d17.ClearFormats()
$ws.Range("E17").Value = "  -1.56%  "

$d18 = $ws.Range("D18")
$d18.NumberFormat = "@"
$d18.Value = "42.598.73"
$d18.ClearFormats()
$ws.Range("E18").Value = "  -1.47%  "

$d19 = $ws.Range("D19")
$d19.NumberFormat = "@"
$d19.Value = "0.0₃0989"
$d19.ClearFormats()
$ws.Range("E19").Value = "  -2.13%  "

$d20 = $ws.Range("D20")
$d20.NumberFormat = "@"
$d20.Value = "6.18"
$d20.ClearFormats()
$ws.Range("E20").Value = "  -3.63%  "

$ws.Range("E21").Value = "  -2.75%  "

$d22 = $ws.Range("D22")
$d22.NumberFormat = "@"
$d22.Value = "232.64"
$d22.ClearFormats()
$ws.Range("E22").Value = "  -1.77%  "

$d23 = $ws.Range("D23")
$d23.NumberFormat = "@"
$d23.Value = "2.16"
$d23.ClearFormats()
$ws.Range("E23").Value = "  -5.87%  "

$ws.Range("E24").Value = "  -3.63%  "

$ws.Range("E25").Value = "  -0.06%  "

$d26 = $ws.Range("D26")
$d26.NumberFormat = "@"
$d26.Value = "11.32"
$d26.ClearFormats()
$ws.Range("E26").Value = "  -4.44%  "

$d27 = $ws.Range("D27")
$d27.NumberFormat = "@"
$d27.Value = "2.32"
$d27.ClearFormats()
$ws.Range("E27").Value = "  -5.38%  "

$d28 = $ws.Range("D28")
$d28.NumberFormat = "@"
$d28.Value = "2.12"
$d28.ClearFormats()
$ws.Range("E28").Value = "  -2.80%  "

$d29 = $ws.Range("D29")
$d29.NumberFormat = "@"
$d29.Value = "170.09"
$d29.ClearFormats()
$ws.Range("E29").Value = "  +1.14%  "

$d30 = $ws.Range("D30")
$d30.NumberFormat = "@"
$d30.Value = "6.71"
$d30.ClearFormats()
$ws.Range("E30").Value = "  +4.99%  "

$d31 = $ws.Range("D31")
$d31.NumberFormat = "@"
$d31.Value = "20.77"
$d31.ClearFormats()
$ws.Range("E31").Value = "  -2.85%  "

$d32 = $ws.Range("D32")
$d32.NumberFormat = "@"
$d32.Value = "0.0853"
$d32.ClearFormats()
$ws.Range("E32").Value = "  +4.56%  "

$ws.Range("E33").Value = "  -6.58%  "

$d34 = $ws.Range("D34")
$d34.NumberFormat = "@"
$d34.Value = "30.25"
$d34.ClearFormats()
$ws.Range("E34").Value = "  -3.48%  "

$d35 = $ws.Range("D35")
$d35.NumberFormat = "@"
$d35.Value = "0.126"
$d35.ClearFormats()
$ws.Range("E35").Value = "  -1.26%  "

$ws.Range("E36").Value = "  -4.64%  "

$ws.Range("E37").Value = "  -2.25%  "

$d38 = $ws.Range("D38")
$d38.NumberFormat = "@"
$d38.Value = "0.0300"
$d38.ClearFormats()
$ws.Range("E38").Value = "  -4.06%  "

$d39 = $ws.Range("D39")
$d39.NumberFormat = "@"
$d39.Value = "13.24"
$d39.ClearFormats()
$ws.Range("E39").Value = "  -3.72%  "

$ws.Range("E40").Value = "  -5.81%  "

$d41 = $ws.Range("D41")
$d41.NumberFormat = "@"
$d41.Value = "5.99"
$d41.ClearFormats()
$ws.Range("E41").Value = "  -2.30%  "

$d42 = $ws.Range("D42")
$d42.NumberFormat = "@"
$d42.Value = "114.63"
$d42.ClearFormats()
$ws.Range("E42").Value = "  +17.69%  "

$ws.Range("E43").Value = "  -3.38%  "

$d44 = $ws.Range("D44")
$d44.NumberFormat = "@"
$d44.Value = "61.05"
$d44.ClearFormats()
$ws.Range("E44").Value = "  -3.36%  "

$d45 = $ws.Range("D45")
$d45.NumberFormat = "@"
$d45.Value = "8.92"
$d45.ClearFormats()
$ws.Range("E45").Value = "  -4.02%  "

$d46 = $ws.Range("D46")
$d46.NumberFormat = "@"
$d46.Value = "0.101"
$d46.ClearFormats()
$ws.Range("E46").Value = "  -3.22%  "

$d47 = $ws.Range("D47")
$d47.NumberFormat = "@"
$d47.Value = "0.999"
$d47.ClearFormats()
$ws.Range("E47").Value = "  -0.26%  "

$d48 = $ws.Range("D48")
$d48.NumberFormat = "@"
$d48.Value = "4.54"
$d48.ClearFormats()
$ws.Range("E48").Value = "  -8.63%  "

$d49 = $ws.Range("D49")
$d49.NumberFormat = "@"
$d49.Value = "1.15"
$d49.ClearFormats()
$ws.Range("E49").Value = "  -4.44%  "

$d50 = $ws.Range("D50")
$d50.NumberFormat = "@"
$d50.Value = "1.17"
$d50.ClearFormats()
$ws.Range("E50").Value = "  -2.84%  "

$d51 = $ws.Range("D51")
$d51.NumberFormat = "@"
$d51.Value = "4.22"
$d51.ClearFormats()
$ws.Range("E51").Value = "  -2.08%  "
